# Apply Arabic translation edits to the document (Email 8&9 template).
$d = $word.ActiveDocument

function Replace-Text($range, $find, $replace) {
    $range.Find.Execute($find, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $replace, 2) | Out-Null
}

# ---------------------------------------------------------------------
# Paragraph 1: "English / Portuguese / French / Thai / Vietnamese / Spanish"
# The first run ("English") sits inside a hyperlink; the COM engine here
# reproduces Word's well known quirk where text typed/replaced right after
# a hyperlink run picks up the hyperlink's blue/underlined character
# formatting. Fix the second run's formatting back to its original
# (red, no underline) once the text has been swapped in.
# ---------------------------------------------------------------------
Replace-Text $d.Content "English" "الإنجليزية"
Replace-Text $d.Content " / Portuguese / French / Thai / Vietnamese / Spanish" " /البرتغالية/الفرنسية/التايلندية/الفيتنامية/الإسبانية"

$fix = $d.Content
$fix.Find.Execute(" /البرتغالية/الفرنسية/التايلندية/الفيتنامية/الإسبانية") | Out-Null
$fix.Font.Color = 255        # 0x0000FF -> RGB FF0000 (red)
$fix.Font.Underline = 0      # no underline

# Paragraph 3: plain "English" heading (own run, unaffected by the quirk).
Replace-Text $d.Content "English" "الإنجليزية"

# Table cell labels.
Replace-Text $d.Content "Brief" "المضمون"
Replace-Text $d.Content "Target audience" "الجمهور المستهدف"
Replace-Text $d.Content "Event attendees" "الحاضرون في الحدث"

# Heading.
Replace-Text $d.Content "We can’t wait to meet you! " "لا يسعنا الانتظار لمقابلتك! "

# Paragraph 15: "Hi [PARTNER NAME], " — scope the Find to this paragraph
# only, since ", " occurs many more times elsewhere in the document.
$greetingPara = $d.Paragraphs.Item(15).Range
Replace-Text $greetingPara "Hi " "مرحبًا  "
$greetingPara2 = $d.Paragraphs.Item(15).Range
Replace-Text $greetingPara2 ", " ",، "

Replace-Text $d.Content "In this email, we’ve linked/attached the following documents:" "في هذه الرسالة الإلكترونية، قمنا بإضافة رابط/إرفاق المستندات التالية:"
Replace-Text $d.Content "Your return flight tickets" "تذاكر رحلة العودة الخاصة بك"
Replace-Text $d.Content "Your accommodation booking details" "تفاصيل حجز الإقامة الخاصة بك"

# Paragraph 24: "If you have any questions, please contact us via live chat or WhatsApp."
Replace-Text $d.Content "If you have any questions, please contact us via " "إذا كانت لديك أي أسئلة، فاتصل بنا:  "

Replace-Text $d.Content "live chat" "الدردشة الحية"
$fix2 = $d.Content
$fix2.Find.Execute("الدردشة الحية") | Out-Null
$fix2.Font.Color = 13391121  # 0xCC5511 -> RGB 1155CC (hyperlink blue)
$fix2.Font.Underline = 1     # single underline

# Paragraph 25.
Replace-Text $d.Content "If you have any questions, please contact your country manager, " "إذا كانت لديك أي أسئلة، فيُرجى الاتصال بمدير بلدك  "

# Comment text lives in the comments part; Find over main Content does not
# reach it, so address it through the comment's own Range.
foreach ($c in $d.Comments) {
    if ($c.Range.Text -eq "choose either one") {
        $c.Range.Text = "اختر أيًا منهما"
    }
}
